# Auto-generated edit script applying odds updates per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.95
$ws.Range("G2").Value = 2.06
$ws.Range("H2").Value = 3.75
$ws.Range("J2").Value = 3.8
$ws.Range("L2").Value = 1.28
$ws.Range("P2").Value = 2.24
$ws.Range("Q2").Value = 1.7
$ws.Range("T2").Value = 1.65
$ws.Range("W2").Value = 1.95
$ws.Range("X2").Value = 24
$ws.Range("Y2").Value = 22
$ws.Range("Z2").Value = 38
$ws.Range("AA2").Value = 90
$ws.Range("AB2").Value = 14
$ws.Range("AC2").Value = 11.5
$ws.Range("AD2").Value = 20
$ws.Range("AE2").Value = 55
$ws.Range("AF2").Value = 17
$ws.Range("AG2").Value = 13
$ws.Range("AH2").Value = 20
$ws.Range("AJ2").Value = 29
$ws.Range("AK2").Value = 24
$ws.Range("AL2").Value = 38
$ws.Range("AM2").Value = 90
$ws.Range("AO2").Value = 44

# Row 3
$ws.Range("F3").Value = 1.6
$ws.Range("G3").Value = 1.68
$ws.Range("H3").Value = 5.7
$ws.Range("L3").Value = 1.27
$ws.Range("M3").Value = 1.04
$ws.Range("Q3").Value = 1.69
$ws.Range("R3").Value = 1.49
$ws.Range("S3").Value = 2.46
$ws.Range("T3").Value = 1.76
$ws.Range("W3").Value = 2.46
$ws.Range("X3").Value = 25
$ws.Range("Y3").Value = 28
$ws.Range("Z3").Value = 60
$ws.Range("AB3").Value = 12
$ws.Range("AC3").Value = 12.5
$ws.Range("AD3").Value = 28
$ws.Range("AE3").Value = 90
$ws.Range("AF3").Value = 13
$ws.Range("AG3").Value = 12.5
$ws.Range("AH3").Value = 24
$ws.Range("AI3").Value = 85
$ws.Range("AJ3").Value = 19.5
$ws.Range("AK3").Value = 20
$ws.Range("AL3").Value = 38
$ws.Range("AM3").Value = 120
$ws.Range("AN3").Value = 9.4

# Row 4
$ws.Range("G4").Value = 1.28
$ws.Range("J4").Value = 5.5

# Row 9
$ws.Range("J9").Value = 5.7
$ws.Range("Q9").Value = 1.65
$ws.Range("Y9").Value = 95
$ws.Range("AA9").Value = 560
$ws.Range("AD9").Value = 46
$ws.Range("AE9").Value = 980
$ws.Range("AH9").Value = 30
$ws.Range("AN9").Value = 5

# Row 10
$ws.Range("H10").Value = 5.8
$ws.Range("I10").Value = 6.2
$ws.Range("J10").Value = 4.8
$ws.Range("Q10").Value = 1.43
$ws.Range("S10").Value = 2.06
$ws.Range("Z10").Value = 1000
$ws.Range("AD10").Value = 24
$ws.Range("AE10").Value = 75
$ws.Range("AF10").Value = 14
$ws.Range("AI10").Value = 1000
$ws.Range("AJ10").Value = 17.5
$ws.Range("AM10").Value = 65
$ws.Range("AN10").Value = 5

# Row 11
$ws.Range("F11").Value = 2.32
$ws.Range("G11").Value = 2.48
$ws.Range("H11").Value = 3.2
$ws.Range("I11").Value = 3.4
$ws.Range("K11").Value = 3.8
$ws.Range("P11").Value = 1.96
$ws.Range("Q11").Value = 1.89

# Row 12
$ws.Range("H12").Value = 2.3
$ws.Range("K12").Value = 3.7

# Row 13
$ws.Range("F13").Value = 1.17
$ws.Range("G13").Value = 2.5
$ws.Range("J13").Value = 2.74
$ws.Range("K13").Value = 980

# Row 14
$ws.Range("F14").Value = 2.26
$ws.Range("G14").Value = 2.28
$ws.Range("H14").Value = 3.4
$ws.Range("I14").Value = 3.45
$ws.Range("N14").Value = 5.1
$ws.Range("P14").Value = 2.32
$ws.Range("Q14").Value = 1.71
$ws.Range("R14").Value = 1.54
$ws.Range("S14").Value = 2.76
$ws.Range("T14").Value = 1.6
$ws.Range("X14").Value = 19.5
$ws.Range("Y14").Value = 17
$ws.Range("AA14").Value = 1000
$ws.Range("AE14").Value = 980
$ws.Range("AI14").Value = 42
$ws.Range("AJ14").Value = 30
$ws.Range("AL14").Value = 32
$ws.Range("AM14").Value = 65
$ws.Range("AN14").Value = 13.5
$ws.Range("AO14").Value = 26

# Row 15
$ws.Range("F15").Value = 3.15
$ws.Range("G15").Value = 3.2
$ws.Range("I15").Value = 2.44
$ws.Range("X15").Value = 21
$ws.Range("Z15").Value = 18.5
$ws.Range("AA15").Value = 34
$ws.Range("AB15").Value = 17
$ws.Range("AC15").Value = 8.800000000000001
$ws.Range("AE15").Value = 23
$ws.Range("AF15").Value = 25
$ws.Range("AH15").Value = 14.5
$ws.Range("AI15").Value = 32
$ws.Range("AJ15").Value = 48
$ws.Range("AK15").Value = 30
$ws.Range("AL15").Value = 44
$ws.Range("AM15").Value = 70

# Row 16
$ws.Range("J16").Value = 5.1
$ws.Range("K16").Value = 5.3
$ws.Range("R16").Value = 1.67
$ws.Range("T16").Value = 1.74
$ws.Range("U16").Value = 2.28
$ws.Range("AA16").Value = 210
$ws.Range("AB16").Value = 11.5
$ws.Range("AC16").Value = 12.5
$ws.Range("AE16").Value = 100
$ws.Range("AF16").Value = 11
$ws.Range("AI16").Value = 75
$ws.Range("AK16").Value = 14.5
$ws.Range("AM16").Value = 110

# Row 17
$ws.Range("I17").Value = 16
$ws.Range("J17").Value = 8
$ws.Range("O17").Value = 1.09
$ws.Range("P17").Value = 3.95
$ws.Range("R17").Value = 2.22
$ws.Range("S17").Value = 1.77
$ws.Range("T17").Value = 1.77
$ws.Range("U17").Value = 2.2
$ws.Range("Y17").Value = 85
$ws.Range("Z17").Value = 210
$ws.Range("AA17").Value = 610
$ws.Range("AB17").Value = 18
$ws.Range("AC17").Value = 21
$ws.Range("AD17").Value = 60
$ws.Range("AF17").Value = 12
$ws.Range("AH17").Value = 30
$ws.Range("AL17").Value = 30

# Row 18
$ws.Range("Q18").Value = 1.44
$ws.Range("U18").Value = 2.04
$ws.Range("X18").Value = 40
$ws.Range("Y18").Value = 1000
$ws.Range("AA18").Value = 510
$ws.Range("AE18").Value = 210
$ws.Range("AF18").Value = 9.6
$ws.Range("AH18").Value = 32
$ws.Range("AM18").Value = 160

# Row 19
$ws.Range("G19").Value = 5.5
$ws.Range("H19").Value = 1.64
$ws.Range("I19").Value = 1.65
$ws.Range("J19").Value = 4.8
$ws.Range("K19").Value = 4.9
$ws.Range("N19").Value = 7.8
$ws.Range("P19").Value = 3.3
$ws.Range("Q19").Value = 1.4
$ws.Range("R19").Value = 1.99
$ws.Range("S19").Value = 1.99
$ws.Range("T19").Value = 1.51
$ws.Range("U19").Value = 2.84
$ws.Range("AA19").Value = 19
$ws.Range("AB19").Value = 36
$ws.Range("AC19").Value = 13
$ws.Range("AE19").Value = 14.5
$ws.Range("AH19").Value = 16
$ws.Range("AJ19").Value = 150
$ws.Range("AM19").Value = 65
$ws.Range("AN19").Value = 34

# Row 20
$ws.Range("F20").Value = 1.53
$ws.Range("G20").Value = 1.54
$ws.Range("H20").Value = 6.8
$ws.Range("I20").Value = 7.2
$ws.Range("J20").Value = 4.8
$ws.Range("O20").Value = 1.21
$ws.Range("Q20").Value = 1.65
$ws.Range("S20").Value = 2.64
$ws.Range("T20").Value = 1.8
$ws.Range("U20").Value = 2.18
$ws.Range("AA20").Value = 210
$ws.Range("AC20").Value = 11.5
$ws.Range("AE20").Value = 95
$ws.Range("AH20").Value = 21
$ws.Range("AJ20").Value = 14.5
$ws.Range("AK20").Value = 15
$ws.Range("AM20").Value = 1000
$ws.Range("AN20").Value = 6.4
$ws.Range("AO20").Value = 1000

# Row 21
$ws.Range("F21").Value = 1.94
$ws.Range("G21").Value = 2.06
$ws.Range("H21").Value = 4.3
$ws.Range("I21").Value = 5.4
$ws.Range("J21").Value = 3.35
$ws.Range("K21").Value = 3.7
$ws.Range("P21").Value = 1.72
$ws.Range("Q21").Value = 2.18

# Row 22
$ws.Range("F22").Value = 1.82
$ws.Range("G22").Value = 1.9
$ws.Range("H22").Value = 5.7
$ws.Range("I22").Value = 6.8
$ws.Range("J22").Value = 3.25
$ws.Range("K22").Value = 3.55
$ws.Range("L22").Value = 1.54
$ws.Range("M22").Value = 1.13
$ws.Range("N22").Value = 2.3
$ws.Range("O22").Value = 1.65
$ws.Range("P22").Value = 1.44
$ws.Range("Q22").Value = 2.98
$ws.Range("S22").Value = 6.6
$ws.Range("T22").Value = 2.46
$ws.Range("U22").Value = 1.53
$ws.Range("V22").Value = 1.17
$ws.Range("W22").Value = 2.1
$ws.Range("X22").Value = 8.199999999999999
$ws.Range("Y22").Value = 14.5
$ws.Range("AA22").Value = 280
$ws.Range("AB22").Value = 5.5
$ws.Range("AC22").Value = 8.800000000000001
$ws.Range("AD22").Value = 1000
$ws.Range("AE22").Value = 190
$ws.Range("AF22").Value = 9.6
$ws.Range("AG22").Value = 12
$ws.Range("AH22").Value = 38
$ws.Range("AI22").Value = 230
$ws.Range("AJ22").Value = 24
$ws.Range("AK22").Value = 32
$ws.Range("AL22").Value = 80
$ws.Range("AM22").Value = 450
$ws.Range("AN22").Value = 1000
